$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- ProductLoan_Input (sheet1) edits ---

# shortname value becomes a plain number (390) instead of the shared string "kar3"
$ws1.Range("B3").Value = 390

# numberofrepaymentsdefault 12 -> 1 -- wait, this is nominalinterestratedefault row (B11)
$ws1.Range("B11").Value = 1

# maximumallowedaoutstandingbalance 5000 -> 10000
$ws1.Range("B28").Value = 10000

# New chart-of-accounts style rows appended below the existing data (31-42).
# Write all of column B first, then column A, so the shared-strings table
# grows in the same order the source workbook used (values before labels).
$ws1.Range("B31").Value = "Cash"
$ws1.Range("B32").Value = "Loan portfolio "
$ws1.Range("B33").Value = "Interest Receivable "
$ws1.Range("B34").Value = "Penalties Receivable "
$ws1.Range("B35").Value = "Transfer in Suspence "
$ws1.Range("B36").Value = "Fees Receivable"
$ws1.Range("B37").Value = "Income from interest"
$ws1.Range("B38").Value = "Income from penalties"
$ws1.Range("B39").Value = "Income from fees"
$ws1.Range("B40").Value = "Income from recovery repayments"
$ws1.Range("B41").Value = "Losses Writtenoff "
$ws1.Range("B42").Value = "Overpayment Liability"

$ws1.Range("A31").Value = "fundsource"
$ws1.Range("A32").Value = "loanprotfolio"
$ws1.Range("A33").Value = "interestreceivable"
$ws1.Range("A34").Value = "penaltiesreceivable"
$ws1.Range("A35").Value = "transferinsuspense"
$ws1.Range("A36").Value = "feesreceivable"
$ws1.Range("A37").Value = "incomefrominterest"
$ws1.Range("A38").Value = "incomefrompenalties"
$ws1.Range("A39").Value = "incomefromfees"
$ws1.Range("A40").Value = "incomefromrecoveryrepayments"
$ws1.Range("A41").Value = "loseswrittenoff"
$ws1.Range("A42").Value = "overpaymentliability"

# Match the formatting used by the rest of column A / column B (fill colours)
# by copying an existing cell's format onto the freshly-added rows.
$ws1.Range("A21").Copy()
$ws1.Range("A31:A42").PasteSpecial(-4122)
$ws1.Range("B10").Copy()
$ws1.Range("B31:B42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen column B to fit the new longer text, and drop the old best-fit flag.
$ws1.Columns.Item(2).ColumnWidth = 58.28

# Scroll/selection bookkeeping + which sheet is active.
$ws1.Activate()
$ws1.Range("B30").Select()

$ws2.Activate()
$ws2.Range("F12").Select()
